# Add "Test Case № 6" (become a member page) to the Shopping Cart Test Suite

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74
$ws.Range("A74").Value = "Test Case № 6"
$ws.Range("B74").Value = "Make sure shopping cart is displayed at the middle of the page"

# Row 75
$ws.Range("A75").Value = "Iteration 1"

# Row 76
$ws.Range("A76").Value = "Steps to reproduce : "
$ws.Range("B76").Value = "1. Navigate to a page containing add to cart button"

# Row 77
$ws.Range("B77").Value = "2. Add some products to the shopping cart"
$ws.Range("C77").Value = "Expected Result : "
$ws.Range("D77").Value = "Products are added to the shopping cart"

# Row 78
$ws.Range("B78").Value = "3. Make sure the shopping cart stays at the middle of the page and products add proper"
$ws.Range("C78").Value = "Expected Result : "
$ws.Range("D78").Value = "Shopping cart is possitioned at the middle of the page"

# Update the selection / view to match author's final cursor position
# (topLeftCell -> A64, active cell -> B66)
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B66").Select()
